# Horarios actualizados Linea 141 - 364
# Updates the "Ultima actualizacion" / "Total filas" headers and the
# schedule data tables on all three worksheets (LP1912, LP1912-215,
# 6203-6173) to reflect the new scrape at 04:18:52.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912  (header row 5, data rows 6-19)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 04:18:52"
$ws1.Range("A3").Value = "Total filas: 14"

$data1 = @(
    @("03:55:13","04:01","81_EL PELIGRO",6,"LP1912"),
    @("03:55:13","04:34","215_ALUAR",39,"LP1912"),
    @("04:18:52","04:45","215A_EL PATO",27,"LP1912"),
    @("03:55:13","04:53","11_ETCHEVERRY",58,"LP1912"),
    @("03:55:13","05:16","17_ROMERO",81,"LP1912"),
    @("04:18:52","05:21","23_HERNANDEZ",63,"LP1912"),
    @("03:55:13","05:22","23_HERNANDEZ",87,"LP1912"),
    @("04:18:52","05:34","215B_EL PATO",76,"LP1912"),
    @("03:55:13","05:35","215B_EL PATO",100,"LP1912"),
    @("03:55:13","05:46","15_ABASTO",111,"LP1912"),
    @("04:18:52","05:53","10_OLMOS",95,"LP1912"),
    @("04:18:52","06:05","16_SANTA ANA",107,"LP1912"),
    @("04:18:52","06:11","215A_EL PATO",113,"LP1912"),
    @("04:18:52","06:13","225_HARAS DEL SUR",115,"LP1912")
)

$r = 6
foreach ($row in $data1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215  (header row 5, data rows 6-10)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 04:18:52"
$ws2.Range("A3").Value = "Total filas: 5"

$data2 = @(
    @("03:55:13","04:34","215_ALUAR",39,"LP1912"),
    @("04:18:52","04:45","215A_EL PATO",27,"LP1912"),
    @("04:18:52","05:34","215B_EL PATO",76,"LP1912"),
    @("03:55:13","05:35","215B_EL PATO",100,"LP1912"),
    @("04:18:52","06:11","215A_EL PATO",113,"LP1912")
)

$r = 6
foreach ($row in $data2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173  (header row 5, data rows 6-8)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 04:18:52"
$ws3.Range("A3").Value = "Total filas: 3"

$data3 = @(
    @("04:18:52","05:43","215A_LA PLATA",85,"L6173"),
    @("03:55:13","05:44","215A_LA PLATA",109,"L6173"),
    @("04:18:52","06:08","215A_LA PLATA",110,"L6173")
)

$r = 6
foreach ($row in $data3) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
